$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph ("Baseline project plan") - apply Calibri font
#    to the whole paragraph (run + paragraph mark).
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Font.Name = "Calibri"
$titlePara.Range.Font.NameBi = "Calibri"

# ------------------------------------------------------------------
# 2) "Team Awesome plans..." paragraph - apply Calibri font to the
#    whole paragraph and change spacing-after from 120 -> 0.
# ------------------------------------------------------------------
$introPara = $d.Paragraphs(3)
$introPara.Range.Font.Name = "Calibri"
$introPara.Range.Font.NameBi = "Calibri"
$introPara.SpaceAfter = 0

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from its old location (just before
#    "manager, no other updates will") to the new location (right
#    after "...keep liquor and beer levels at the "). Word only
#    allows one bookmark with a given name, so adding it at the new
#    spot automatically removes it from the old spot - this also
#    splits the run exactly like the authentic edit did.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("keep liquor and beer levels at the ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPos = $find.Parent.End
    $bmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ------------------------------------------------------------------
# 4) Split "Paul Naumann" into two runs ("Paul " / "Naumann") the
#    way Word's spell-checker does when it flags "Naumann". We force
#    a run boundary at that position by toggling a character
#    attribute on/off (Word does not re-merge the runs afterwards).
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("Naumann", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $naumannRange = $find2.Parent
    $naumannRange.Font.Bold = $true
    $naumannRange.Font.Bold = $false
}

Write-Host "done"
